$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing gpt-4o / gpt-4o-2024-05-13 down by one.
$ws.Rows.Item(2).Insert()

# Copy the style (center alignment) from the row below into the newly inserted row.
$ws.Range("A2:AA2").HorizontalAlignment = -4108
$ws.Range("A2:AA2").VerticalAlignment = -4108

# Fill in the new row 2 with the gpt-4.1 results.
$rowData = @(
    "gpt-4.1", 0.9399999999999999, 0.98, 1, 0.9, 0.9, 1, 0.953, 0.91, 0.89, 0.9,
    0.64, 0.64, 0.64, 0.76, 0.8, 0.863, 0.9, 0.9399999999999999, 0.88, 0.907,
    0.5669999999999999, 0.746, 0.95, 0.958, 0.758, 0.859
)

for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $rowData[$i]
}
